# Refresh the cryptos price list (GitHub Actions scheduled update).
# Prices/volumes in column D/E are stored as *text* (European dotted
# grouping like "42.976.89", or fixed-decimal strings like "0.810"/"10.10"
# that must keep trailing zeros), so each D-cell is pre-formatted as Text
# ("@") before the value is assigned - otherwise Excel's smart-entry would
# coerce the numeric-looking string into a real number/date and mangle it
# (dropping trailing zeros, misreading the dotted grouping, etc).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.976.89"
$ws.Range("E2").Value = "  +4.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.245.63"
$ws.Range("E3").Value = "  +3.23%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.39"
$ws.Range("E5").Value = "  +3.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.17"
$ws.Range("E7").Value = "  +8.26%  "

$ws.Range("E8").Value = "  -0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.619"
$ws.Range("E9").Value = "  +6.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.13"
$ws.Range("E10").Value = "  +1.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0935"
$ws.Range("E11").Value = "  +0.84%  "

$ws.Range("E12").Value = "  +2.67%  "

$ws.Range("E13").Value = "  +1.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.562.17"
$ws.Range("E14").Value = "  +2.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.63"
$ws.Range("E15").Value = "  +4.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.243.33"
$ws.Range("E16").Value = "  +3.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.810"
$ws.Range("E17").Value = "  +0.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.881.34"
$ws.Range("E18").Value = "  +4.47%  "

$ws.Range("E19").Value = "  +3.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.23"
$ws.Range("E20").Value = "  +1.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.01"
$ws.Range("E21").Value = "  +0.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.10"
$ws.Range("E22").Value = "  +2.46%  "

# Row 23: coin swapped places with its neighbour in the ranking
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.37"
$ws.Range("E23").Value = "  +2.26%  "

# Row 24: coin swapped places with its neighbour in the ranking
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.21"
$ws.Range("E24").Value = "  +13.71%  "

$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.90"
$ws.Range("E26").Value = "  -0.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.43"
$ws.Range("E27").Value = "  -3.46%  "

$ws.Range("E28").Value = "  +2.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.43"
$ws.Range("E29").Value = "  +22.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "173.64"
$ws.Range("E30").Value = "  +3.56%  "

$ws.Range("E31").Value = "  -1.73%  "

$ws.Range("E32").Value = "  +1.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0799"
$ws.Range("E33").Value = "  +3.69%  "

$ws.Range("E34").Value = "  +3.37%  "

$ws.Range("E35").Value = "  +1.23%  "

$ws.Range("E36").Value = "  +6.47%  "

$ws.Range("E37").Value = "  +4.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0331"
$ws.Range("E38").Value = "  +15.67%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.84"
$ws.Range("E39").Value = "  +8.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.15"
$ws.Range("E40").Value = "  +3.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.56"
$ws.Range("E41").Value = "  +2.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.204"
$ws.Range("E42").Value = "  +6.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "107.53"
$ws.Range("E43").Value = "  +9.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "59.82"
$ws.Range("E44").Value = "  -0.66%  "

$ws.Range("E45").Value = "  +4.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0989"
$ws.Range("E46").Value = "  +1.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.469"
$ws.Range("E47").Value = "  +23.29%  "

$ws.Range("E48").Value = "  +5.67%  "

$ws.Range("E49").Value = "  +1.59%  "

$ws.Range("E50").Value = "  +1.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.449.81"
